$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model_Matched_Parameters")
$ws.Range("K115").Value = 0.0033
Write-Output "done"
